# Replace the double/curly quotation marks that wrap the spoken dialogue
# with straight single quotes, leaving the [name="..."] tag quotes intact.
# (commit: "update on 20210731 画中人")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C41").Value  = '[name="Frostleaf"]   ''Wherever there is blood...''' + "`n"
$ws.Range("C53").Value  = '[name="Amiya"]   ''...Your mind can hide nothing from me.''' + "`n"
$ws.Range("C97").Value  = '[name="FrostNova"]   ''Sleep, my children. Sleep.' + "`n"
$ws.Range("C98").Value  = '[name="FrostNova"]   Hedgehog dolls and teddy bears...''......' + "`n"
$ws.Range("C101").Value = '[name="FrostNova"]   ''Sink quietly into black...♪' + "`n"
$ws.Range("C102").Value = '[name="FrostNova"]   ...My broken puppets and dolls♬''' + "`n"
